$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Move the "test number" generator up to row 1 (D1/E1) ---
$ws.Range("D1").Formula = "=SUM(C:C)"
$ws.Range("E1").Formula = "=RANDBETWEEN(1000,9999)"

# Remove the old D2/E2 cells (their content now lives in D1/E1)
$ws.Range("D2:E2").ClearContents()

# --- Repoint the IF() formulas in column C (rows 2-16) at $E$1 instead of $E$2 ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("C$r").Formula = "=IF(B$r=`$E`$1,1,0)"
}

# --- Second "лр2" block: rows 17-56 ---

# Copy the "лр1" block's formatting onto the new A17:A56 range so the
# merged label cell picks up the same centred style (s="1") instead of
# Excel's unstyled default.
$ws.Range("A3").Copy()
$ws.Range("A17:A56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A17").Value = "лр2"

$b17 = @(4411,9298,3354,8569,3730,9007,5201,2981,4312,8428,2361,5063,7711,8833,1262,9020,1934,7237,3943,7619,3832,1346,9622,8873,7799,9354,9130,5895,2461,2790,2624,5871,9164,7457,9865,3591,3558,4366,5789,6522)

for ($i = 0; $i -lt $b17.Length; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = $b17[$i]
    $ws.Range("C$r").Formula = "=IF(B$r=`$E`$1,1,0)"
}

$ws.Range("A17:A56").Merge()

# --- View: freeze the header row, then select/scroll to the new block ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A17:A56").Select()
$excel.ActiveWindow.ScrollRow = 22
